# Auto-generated edit script applying numeric corrections to the
# profit-calculation sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1165716.6
$ws.Range("J17").Value = 1165716.6
$ws.Range("L17").Value = 3497149.8
$ws.Range("N17").Value = -3497485.8

# Row 87
$ws.Range("H87").Value = 26779.023
$ws.Range("J87").Value = 26779.023
$ws.Range("L87").Value = 26779.023
$ws.Range("N87").Value = -29275.023

# Row 90
$ws.Range("H90").Value = 26779.023
$ws.Range("J90").Value = 26779.023
$ws.Range("L90").Value = 80337.069
$ws.Range("N90").Value = -92817.069

# Row 101
$ws.Range("H101").Value = 47619412
$ws.Range("I101").Value = 47619412
$ws.Range("K101").Value = 142858236
$ws.Range("M101").Value = -142856614

# Row 113
$ws.Range("H113").Value = 10337.526
$ws.Range("I113").Value = 4124.75
$ws.Range("J113").Value = 11994.267
$ws.Range("K113").Value = 4124.75
$ws.Range("L113").Value = 11994.267
$ws.Range("M113").Value = -870.75
$ws.Range("N113").Value = -18502.267

# Row 138
$ws.Range("H138").Value = 3189.7866
$ws.Range("I138").Value = 1561.7435
$ws.Range("K138").Value = 4685.2305
$ws.Range("M138").Value = 454.7695000000003

# Row 140
$ws.Range("H140").Value = 62641
$ws.Range("J140").Value = 62641
$ws.Range("L140").Value = 62641
$ws.Range("N140").Value = -73001


# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1648
$ws.Range("I61").Value = 1081.7
$ws.Range("J61").Value = 2355.875
$ws.Range("K61").Value = 1081.7
$ws.Range("L61").Value = 2355.875
$ws.Range("M61").Value = -869.7
$ws.Range("N61").Value = -2779.875

# Row 74
$ws.Range("H74").Value = 2543.0938
$ws.Range("I74").Value = 2782.5
$ws.Range("J74").Value = 867.25
$ws.Range("K74").Value = 2782.5
$ws.Range("L74").Value = 867.25
$ws.Range("M74").Value = -1908.5
$ws.Range("N74").Value = -2615.25

# Row 77
$ws.Range("H77").Value = 2543.0938
$ws.Range("I77").Value = 2782.5
$ws.Range("J77").Value = 867.25
$ws.Range("K77").Value = 13912.5
$ws.Range("L77").Value = 4336.25
$ws.Range("M77").Value = -9544.5
$ws.Range("N77").Value = -13072.25

# Row 94
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51802

# Row 132
$ws.Range("H132").Value = 1985.9048
$ws.Range("I132").Value = 1056.32
$ws.Range("J132").Value = 3352.9412
$ws.Range("K132").Value = 3168.96
$ws.Range("L132").Value = 10058.8236
$ws.Range("M132").Value = -638.96
$ws.Range("N132").Value = -15118.8236

# Row 136
$ws.Range("H136").Value = 1648
$ws.Range("I136").Value = 1081.7
$ws.Range("J136").Value = 2355.875
$ws.Range("K136").Value = 3245.1
$ws.Range("L136").Value = 7067.625
$ws.Range("M136").Value = -695.1000000000004
$ws.Range("N136").Value = -12167.625


# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1265.2222
$ws.Range("I94").Value = 1517.4
$ws.Range("J94").Value = 950
$ws.Range("K94").Value = 1517.4
$ws.Range("L94").Value = 950
$ws.Range("M94").Value = -1066.4
$ws.Range("N94").Value = -1852

# Row 103
$ws.Range("H103").Value = 19333.334
$ws.Range("J103").Value = 19333.334
$ws.Range("L103").Value = 19333.334
$ws.Range("N103").Value = -21677.334


# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2330.1804
$ws.Range("I31").Value = 1591.4117
$ws.Range("J31").Value = 3260.4814
$ws.Range("K31").Value = 1591.4117
$ws.Range("L31").Value = 3260.4814
$ws.Range("M31").Value = -1296.4117
$ws.Range("N31").Value = -3850.4814

# Row 34
$ws.Range("H34").Value = 2330.1804
$ws.Range("I34").Value = 1591.4117
$ws.Range("J34").Value = 3260.4814
$ws.Range("K34").Value = 1591.4117
$ws.Range("L34").Value = 3260.4814
$ws.Range("M34").Value = -1389.4117
$ws.Range("N34").Value = -3664.4814

# Row 58
$ws.Range("H58").Value = 1621.34
$ws.Range("I58").Value = 1198.1765
$ws.Range("J58").Value = 2520.5625
$ws.Range("K58").Value = 1198.1765
$ws.Range("L58").Value = 2520.5625
$ws.Range("M58").Value = -995.1765
$ws.Range("N58").Value = -2926.5625

# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("N87").ClearContents()

# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("N90").ClearContents()

# Row 107
$ws.Range("H107").Value = 894.4
$ws.Range("I107").Value = 594.9167
$ws.Range("J107").Value = 1343.625
$ws.Range("K107").Value = 594.9167
$ws.Range("L107").Value = 1343.625
$ws.Range("M107").Value = 1325.0833
$ws.Range("N107").Value = -5183.625

# Row 134
$ws.Range("H134").Value = 1770.9788
$ws.Range("I134").Value = 1590.439
$ws.Range("J134").Value = 3004.6667
$ws.Range("K134").Value = 4771.317
$ws.Range("L134").Value = 9014.000100000001
$ws.Range("M134").Value = -2236.317
$ws.Range("N134").Value = -14084.0001

# Row 136
$ws.Range("H136").Value = 1621.34
$ws.Range("I136").Value = 1198.1765
$ws.Range("J136").Value = 2520.5625
$ws.Range("K136").Value = 3594.5295
$ws.Range("L136").Value = 7561.6875
$ws.Range("M136").Value = -1044.5295
$ws.Range("N136").Value = -12661.6875


# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1227.2727
$ws.Range("I68").Value = 460
$ws.Range("J68").Value = 1866.6666
$ws.Range("K68").Value = 1380
$ws.Range("L68").Value = 5599.9998
$ws.Range("M68").Value = -569
$ws.Range("N68").Value = -7221.9998

# Row 71
$ws.Range("H71").Value = 1227.2727
$ws.Range("I71").Value = 460
$ws.Range("J71").Value = 1866.6666
$ws.Range("K71").Value = 4140
$ws.Range("L71").Value = 16799.9994
$ws.Range("M71").Value = -84
$ws.Range("N71").Value = -24911.9994

# Row 107
$ws.Range("H107").Value = 1040723.06
$ws.Range("J107").Value = 1689630.9
$ws.Range("L107").Value = 5068892.699999999
$ws.Range("N107").Value = -5072732.699999999

# Row 131
$ws.Range("H131").Value = 6396.6
$ws.Range("J131").Value = 9387.154
$ws.Range("L131").Value = 28161.462
$ws.Range("N131").Value = -38241.462

# Row 138
$ws.Range("H138").Value = 1508.4584
$ws.Range("I138").Value = 1022.5
$ws.Range("K138").Value = 3067.5
$ws.Range("M138").Value = 2072.5

# Row 139
$ws.Range("H139").Value = 1491.1111
$ws.Range("I139").Value = 1017.3333
$ws.Range("J139").Value = 2083.3333
$ws.Range("K139").Value = 3051.9999
$ws.Range("L139").Value = 6249.999899999999
$ws.Range("M139").Value = 2088.0001
$ws.Range("N139").Value = -16529.9999

# Row 141
$ws.Range("H141").Value = 3396.6667
$ws.Range("I141").Value = 752
$ws.Range("J141").Value = 5285.7144
$ws.Range("K141").Value = 2256
$ws.Range("L141").Value = 15857.1432
$ws.Range("M141").Value = 2924
$ws.Range("N141").Value = -26217.1432


# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 957.8182
$ws.Range("I100").Value = 964.6
$ws.Range("J100").Value = 890
$ws.Range("K100").Value = 964.6
$ws.Range("L100").Value = 890
$ws.Range("M100").Value = -423.6
$ws.Range("N100").Value = -1972

# Row 136
$ws.Range("H136").Value = 8773547
$ws.Range("I136").Value = 1674.3235
$ws.Range("J136").Value = 83334460
$ws.Range("K136").Value = 5022.970499999999
$ws.Range("L136").Value = 250003380
$ws.Range("M136").Value = -2472.970499999999
$ws.Range("N136").Value = -250008480


# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 1624.2667
$ws.Range("I81").Value = 996
$ws.Range("J81").Value = 1938.4
$ws.Range("K81").Value = 1992
$ws.Range("L81").Value = 3876.8
$ws.Range("M81").Value = -931
$ws.Range("N81").Value = -5998.8

# Row 84
$ws.Range("H84").Value = 1624.2667
$ws.Range("I84").Value = 996
$ws.Range("J84").Value = 1938.4
$ws.Range("K84").Value = 9960
$ws.Range("L84").Value = 19384
$ws.Range("M84").Value = -4656
$ws.Range("N84").Value = -29992

# Row 132
$ws.Range("H132").Value = 1104.6897
$ws.Range("I132").Value = 901.3205
$ws.Range("J132").Value = 2867.2222
$ws.Range("K132").Value = 2703.9615
$ws.Range("L132").Value = 8601.6666
$ws.Range("M132").Value = -173.9615000000003
$ws.Range("N132").Value = -13661.6666

# Row 136
$ws.Range("H136").Value = 5377845
$ws.Range("I136").Value = 7576290.5
$ws.Range("J136").Value = 3866.6667
$ws.Range("K136").Value = 22728871.5
$ws.Range("L136").Value = 11600.0001
$ws.Range("M136").Value = -22726321.5
$ws.Range("N136").Value = -16700.0001

